# Apply the "log" update for the latest dev-log entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: fill in the task description & hours for 2024-11-28 (value already present in A10)
$ws.Range("B10").Value = "item inventory done, and more front end for both dialogue and te rest of the inventory"
$ws.Range("C10").Value = 5

# Row 11: new date entry for 2024-11-29 (same date style as the rows above it)
$ws.Range("A11").Value = 45625
$ws.Range("A11").NumberFormat = $ws.Range("A10").NumberFormat

# Update selection to reflect where the user left off editing
$ws.Range("B11").Select()
